# Reduce the "space before" setting on the first four paragraphs of the
# code-block content placeholder on slide 18 from 3pt to 1pt.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

for ($i = 1; $i -le 4; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.ParagraphFormat.SpaceBefore = 1
}
